# Financials Yearly update - insert new reporting period column before column D
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift existing data (D7:D102 -> E7:L102) to make room for the new period column D
$ws.Range("D7:D102").Insert(-4161)   # xlShiftToRight

# Copy number formats (date / number styles) from the now-shifted column E into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Match the new column's width to its neighbour
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the latest reporting-period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 575900
$ws.Range("D9").Value = 412500
$ws.Range("D10").Value = 163400
$ws.Range("D12").Value = 3600
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 8900
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 543100
$ws.Range("D18").Value = 32800
$ws.Range("D20").Value = 500
$ws.Range("D21").Value = 61900
$ws.Range("D22").Value = 4200
$ws.Range("D23").Value = 29100
$ws.Range("D24").Value = 4600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 24600
$ws.Range("D27").Value = 24400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("D33").Value = 24400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 24400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 93600
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 93900
$ws.Range("D44").Value = 158500
$ws.Range("D45").Value = 4100
$ws.Range("D46").Value = 350100
$ws.Range("D47").Value = 5400
$ws.Range("D48").Value = 213300
$ws.Range("D49").Value = 35300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 12800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 616900
$ws.Range("D57").Value = 55800
$ws.Range("D58").Value = 10500
$ws.Range("D59").Value = 45000
$ws.Range("D60").Value = 111300
$ws.Range("D61").Value = 7100
$ws.Range("D62").Value = 32100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 144800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 360700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 472100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 24400
$ws.Range("D83").Value = 28600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 14700
$ws.Range("D91").Value = -21000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -21200
$ws.Range("D96").Value = -20300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -38400
$ws.Range("D101").Value = -300
$ws.Range("D102").Value = -45100

"Done"
